# "Colocando header nos graficos" - add a header label to column A on each
# data sheet, drop the bold/bordered header style from the row-label cells
# (they keep plain formatting now that the real header lives in row 1), and
# fix a handful of accented Portuguese labels that had lost their diacritics.
# Sheet5 also loses its now-unused "Teto" row, and sheet6 gets a real "2015"
# column header plus refreshed totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share an identical layout: row 1 = year headers (B:E), rows
# 2-12 = source/technology rows. Add the "Fonte/Tecnologia" header to A1,
# strip the header style from A2:A12, and fix a few labels that were
# missing their accents.
# ---------------------------------------------------------------------
$sourceSheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $sourceSheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New header cell A1 - clone the existing header style from B1 so it
    # matches the other header cells exactly, then set its text.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # The row labels (A2:A12) no longer carry the bold/bordered header
    # style - only the actual header row does now.
    $ws.Range("A2:A12").Style = "Normal"

    # Restore missing diacritics.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."
}

# ---------------------------------------------------------------------
# Sheet 5 ("Emissoes Totais (MtCO2eq)"): add "Período" header, fix labels,
# and remove the unused "Teto" row (row 4).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2:A3").Style = "Normal"
$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Range("A4:E4").EntireRow.Delete()

# ---------------------------------------------------------------------
# Sheet 6 ("Custo Total (bilhões de R$)"): add "Tipo Expansão" header,
# turn the old "Custo" header into a "2015" year header (copied from
# another sheet so it stays text, matching the rest of the workbook),
# fix labels, and update the two totals.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsYear = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# Copy a real "2015" text header (value + style) so B1 stays text, not a
# number, exactly like every other sheet's year headers.
$wsYear.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4104)

$ws6.Range("A2:A3").Style = "Normal"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 626

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
